$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header in A1 from "First Word" to "Car Names"
$ws.Range("A1").Value = "Car Names"

# Set text for J1 (was an empty inline string cell)
$ws.Range("J1").Value = "Predicted headform score (excluding blue points)"
